$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet after Problem1 and name it Problem2
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Problem2"

# ---- Values ----
# Header row
$ws2.Range("A1").Value = "Benchmark"
$ws2.Range("B1").Value = "Basic Blocks"
$ws2.Range("C1").Value = "Memory Reads"
$ws2.Range("D1").Value = "Memory Writes"
$ws2.Range("E1").Value = "Total Instructions"

# Row 2: mm_mult_serial 256x256
$ws2.Range("A2").Value = "mm_mult_serial 256x256"
$ws2.Range("B2").Value = 2431
$ws2.Range("C2").Value = 3816
$ws2.Range("D2").Value = 2559
$ws2.Range("E2").Value = 482639275

# Row 3: accumulate array 1x256
$ws2.Range("A3").Value = "accumulate array 1x256"
$ws2.Range("B3").Value = 3429
$ws2.Range("C3").Value = 5634
$ws2.Range("D3").Value = 4632
$ws2.Range("E3").Value = 1421207

# Row 4: 621_wrf_s train
$ws2.Range("A4").Value = "621_wrf_s train"
$ws2.Range("B4").Value = 61000
$ws2.Range("C4").Value = 260281
$ws2.Range("D4").Value = 251764
$ws2.Range("E4").Value = 583459817248

# Row 5: 623_xalanchbmk_s train
$ws2.Range("A5").Value = "623_xalanchbmk_s train"
$ws2.Range("B5").Value = 23105
$ws2.Range("C5").Value = 53870
$ws2.Range("D5").Value = 39434
$ws2.Range("E5").Value = 257156226792

# ---- Formatting: reuse existing styles from Problem1 via copy/paste-special ----
# Header style (yellow fill + border) from Problem1!A1
$ws1.Range("A1").Copy()
$ws2.Range("A1:E1").PasteSpecial(-4122)

# Text-body style (border only) from Problem1!A2
$ws1.Range("A2").Copy()
$ws2.Range("A2:A5").PasteSpecial(-4122)

# Number-body style (border + 2-decimal number format) from Problem1!B2
$ws1.Range("B2").Copy()
$ws2.Range("B2:E5").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Column widths (approximate best-fit widths matching Problem1's look & feel)
$ws2.Columns.Item(1).ColumnWidth = 22.42515625
$ws2.Columns.Item(2).ColumnWidth = 10.58578125
$ws2.Columns.Item(3).ColumnWidth = 13.42515625
$ws2.Columns.Item(4).ColumnWidth = 13.92546875
$ws2.Columns.Item(5).ColumnWidth = 15.58578125

# Page setup: portrait orientation to match Problem1
$ws2.PageSetup.Orientation = 1

# ---- Selections ----
# Update selection on Problem1 sheet first (it is not the active tab afterwards)
$ws1.Range("A2:A3").Select()

# Finally activate Problem2 and set its selection so it becomes the active/tabSelected sheet
$ws2.Activate()
$ws2.Range("H13").Select()
